$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy header style from an existing header cell (e.g. AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Fill in team record data for each data row (2 through 52)
for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = 64   # AD
    $ws.Cells.Item($r, 31).Value = 98   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
